$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "2022-08-30 17:08:31"
$ws.Range("A4").Value = "2022-08-30 17:08:32"
$ws.Range("A5").Value = "2022-08-30 17:08:33"
$ws.Range("D5").Value = "Buscar valores para cada indicador: DSD TX NEW"
$ws.Range("A6").Value = "2022-08-30 17:08:38"
$ws.Range("D6").Value = "Buscar valores para cada indicador: DSD TX CURR"
